# Generate Report for Handback
# Update timestamp strings recorded for the fe9c4b10-981f-41cf-92f9-11aa7b0ec7d4 file
# in the Overview, zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for fe9c4b10-... row
$wsOverview.Range("G4").Value = "2016-08-24 08:47:11"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for fe9c4b10-... row
$wsZhCn.Range("H4").Value = "2016-08-24 08:46:58"
$wsZhCn.Range("K4").Value = "2016-08-24 08:47:30"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for fe9c4b10-... row
$wsDeDe.Range("H4").Value = "2016-08-24 08:47:11"
$wsDeDe.Range("K4").Value = "2016-08-24 08:47:37"
